$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "time_taken" header in F1, matching the style of the other headers (E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Populate F2:F17 with the recorded time_taken values (as plain text)
$ws.Range("F2").Value = "2021-10-05 13:40:49.566809"
$ws.Range("F3").Value = "2021-10-05 13:40:49.566820"
$ws.Range("F4").Value = "2021-10-05 13:40:49.566823"
$ws.Range("F5").Value = "2021-10-05 13:40:49.566826"
$ws.Range("F6").Value = "2021-10-05 13:40:49.566829"
$ws.Range("F7").Value = "2021-10-05 13:40:49.566831"
$ws.Range("F8").Value = "2021-10-05 13:40:49.566834"
$ws.Range("F9").Value = "2021-10-05 13:40:49.566836"
$ws.Range("F10").Value = "2021-10-05 13:40:49.566839"
$ws.Range("F11").Value = "2021-10-05 13:40:49.566842"
$ws.Range("F12").Value = "2021-10-05 13:40:49.566844"
$ws.Range("F13").Value = "2021-10-05 13:40:49.566847"
$ws.Range("F14").Value = "2021-10-05 13:40:49.566849"
$ws.Range("F15").Value = "2021-10-05 13:40:49.566852"
$ws.Range("F16").Value = "2021-10-05 13:40:49.566854"
$ws.Range("F17").Value = "2021-10-05 13:40:49.566857"
